$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the context/cue columns (C,D) for rows 8 and 9 - part of
# randomizing the context roles instead of hard-coding them.
$ws.Range("C8:D9").ClearContents()

# Move the active selection to D9
$ws.Range("D9").Select()
